$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.525.86'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -1.26%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.514.04'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  +0.00%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.16'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.69%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.42'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  +0.59%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.512.63'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("E10").Value = '  -1.39%  '

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.167'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("E12").Value = '  +4.34%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.89'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +0.91%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.973.59'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -0.29%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.329.68'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -1.35%  '

$ws.Range("E16").Value = '  -1.78%  '

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.80'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -1.07%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.516.63'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("E19").Value = '  -1.67%  '

$ws.Range("E20").Value = '  -0.82%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.01'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -2.43%  '

$ws.Range("E22").Value = '  -1.40%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.00'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.92%  '

$ws.Range("E24").Value = '  -0.01%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.23'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("E26").Value = '  -3.20%  '

$ws.Range("E27").Value = '  -1.83%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.640.79'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("E29").Value = '  +0.76%  '

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0890'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -2.54%  '

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.81'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -0.56%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '463.15'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -4.45%  '

$ws.Range("E33").Value = '  -5.01%  '

$ws.Range("E34").Value = '  -2.03%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("E36").Value = '  +1.89%  '

$ws.Range("E37").Value = '  +0.97%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.08'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -0.91%  '

$ws.Range("E42").Value = '  -1.81%  '

$ws.Range("E43").Value = '  -2.85%  '

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.08'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -0.63%  '

$ws.Range("E45").Value = '  -7.60%  '

$ws.Range("E46").Value = '  -7.72%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.46'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -1.55%  '

$ws.Range("E48").Value = '  -1.54%  '

$ws.Range("E49").Value = '  -2.03%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0733'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("E51").Value = '  -3.58%  '
